$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing header texts for columns C..L (3..12) before shifting them,
# since writes below would otherwise clobber values we still need to read.
$captured = @{}
for ($c = 3; $c -le 12; $c++) {
    $captured[$c] = $ws.Cells.Item(1, $c).Text
}

# Shift old C..K (3..11) one column to the right, into D..L (4..12),
# working right-to-left so we never overwrite a cell before it is copied.
for ($c = 11; $c -ge 3; $c--) {
    $ws.Cells.Item(1, $c + 1).Value = $captured[$c]
}

# New column C gets the new "Nº Checkout" header.
$ws.Cells.Item(1, 3).Value = "Nº Checkout"

# The old L1 ("E-mail do Comprador") moves to the end of the header row (AH1).
$ws.Cells.Item(1, 34).Value = $captured[12]

# Give the new column a sensible best-fit-style width (closest the engine can represent).
$ws.Columns(3).ColumnWidth = 10.6
